# Edit script: renumber functional requirements (RF) list.
#
# The original "Requisiti funzionali" section had 14 items (RF1..RF14).
# The old RF1 paragraph ("Se il magazzino offre...") is removed entirely,
# and every following item's number shifts down by one (old RF2 becomes
# the new RF1, old RF3 becomes the new RF2, etc.), plus a handful of
# wording tweaks on specific items.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the old "RF1: Se il magazzino offre..." paragraph
# entirely (its whole paragraph, including the paragraph mark). This
# is the paragraph right after the "Requisiti funzionali" heading and
# the blank paragraph that follows it.
# ---------------------------------------------------------------------

# Locate the paragraph that starts with "RF1: Se il magazzino offre"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "RF1: Se il magazzino offre*") {
        $target = $p
        break
    }
}
$target.Range.Delete()

# ---------------------------------------------------------------------
# Step 2: renumber the remaining items (old RF2..RF14 -> new RF1..RF13)
# and apply the small wording changes from the diff.
# ---------------------------------------------------------------------

# old RF2 -> RF1 : "creare un ordine per un tavolo"
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF2: Il sistema deve fornire al cameriere una funzionalità per creare un ordine per un tavolo*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF2:", $false, $false, $false, $false, $false, $true, 1, $false, "RF1:", 2)

# old RF3 -> RF2 : "inserire pietanze singole"
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF3: Il sistema deve fornire al cameriere una funzionalità per inserire pietanze*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF3:", $false, $false, $false, $false, $false, $true, 1, $false, "RF2:", 2)

# old RF4 -> RF3 : "inserire menu fissi"
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF4: Il sistema deve fornire al cameriere una funzionalità per inserire menu fissi*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF4:", $false, $false, $false, $false, $false, $true, 1, $false, "RF3:", 2)

# old RF5 -> RF4 : "confermare l'ordine"
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF5: Il sistema deve fornire al cameriere una funzionalità per confermare*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF5:", $false, $false, $false, $false, $false, $true, 1, $false, "RF4:", 2)

# old RF6 -> RF5 : "inviare gli ordini validi" + add trailing comment sentence
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF6: Il sistema deve fornire al cameriere una funzionalità che permette di inviare*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF6:", $false, $false, $false, $false, $false, $true, 1, $false, "RF5:", 2)
$endOfPara = $p.Range
$endOfPara.SetRange($endOfPara.End - 1, $endOfPara.End - 1)
$endOfPara.InsertAfter(" -direi di toglierlo perché è un insieme di azioni che fa parte della conferma degli ordini-")

# old RF7 -> RF6 : "registrare il numero di posti"
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF7: Il sistema deve fornire al cameriere un modo per registrare*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF7:", $false, $false, $false, $false, $false, $true, 1, $false, "RF6:", 2)

# old RF8 -> RF7 : "effettuare un controllo" -> "poter effettuare  controlli"
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF8: Il sistema deve effettuare un controllo sul magazzino*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF8: Il sistema deve effettuare un controllo sul magazzino", $false, $false, $false, $false, $false, $true, 1, $false, "RF7: Il sistema deve poter effettuare  controlli sul magazzino", 2)

# old RF9 -> RF8 : "dopo aver confermato un ordine...prenotare"
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF9: Il sistema, dopo aver confermato*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF9:", $false, $false, $false, $false, $false, $true, 1, $false, "RF8:", 2)

# old RF10 -> RF9 : "IL sistema deve presentare..."
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF10: IL sistema deve presentare*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF10:", $false, $false, $false, $false, $false, $true, 1, $false, "RF9:", 2)

# old RF11 -> RF10 : "offrire al cuoco...elenco degli ordini prelevati"
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF11: Il sistema deve offrire al cuoco*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF11:", $false, $false, $false, $false, $false, $true, 1, $false, "RF10:", 2)

# old RF12 -> RF11 : "fornire al cuoco...prelevare il prossimo ordine"
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF12: Il sistema deve fornire al cuoco*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF12:", $false, $false, $false, $false, $false, $true, 1, $false, "RF11:", 2)

# old RF13 -> RF12 : "fornire al direttore..."
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF13: Il sistema deve fornire al direttore*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF13:", $false, $false, $false, $false, $false, $true, 1, $false, "RF12:", 2)

# old RF14 -> RF13 : "fornire al cassiere..." + trailing period added
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "RF14: Il sistema deve fornire al cassiere*") {
        $p = $para
        break
    }
}
$p.Range.Find.Execute("RF14:", $false, $false, $false, $false, $false, $true, 1, $false, "RF13:", 2)
$endOfPara = $p.Range
$endOfPara.SetRange($endOfPara.End - 1, $endOfPara.End - 1)
$endOfPara.InsertAfter(".")
